$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Fix the two comment texts whose wording changed in Sheet1 ---
$ws1.Range("A98").Value = "Đt Sài cũng ổn mỗi cái hay Ci Zalo gọi video call tự tắt cam r không kết thúc cuộc gọi được"
$ws1.Range("A167").Value = "Messenger không nhận được thông báo cuộc gọi và tin nhắn lỗi tùm lum không hài lòng"

# --- Add the two new sheets, right after Sheet1 ---
$wsTot = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$wsTot.Name = "Tốt"

$wsKhongTot = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsTot)
$wsKhongTot.Name = "Không tốt"

# --- Populate "Tốt" sheet (positive sentiment words, sorted) ---
$totWords = @(
 " chắc chắn",
 " cực trâu",
 " Khá là ổn",
 " khá mượt",
 " khá ok",
 " khá ổn",
 " mượt",
 "bền",
 "bền",
 "chất lượng",
 "chất lượng",
 "cực nhanh",
 "dễ sữ dụng",
 "đẹp",
 "đẹp",
 "dùng ổn",
 "Dùng rất tốt",
 "dùng tốt",
 "giá cả hợp lý",
 "giới thiệu thêm",
 "hài lòng",
 "hay",
 "hoàn hảo",
 "hợp lí",
 "hợp lý",
 "khá",
 "khỏi phải chê",
 "không bị bong",
 "ko mơ gì hơn",
 "mạnh",
 "màu sắc đẹp",
 "máy mượt",
 "mượt",
 "mượt mà",
 "nên mua",
 "nên tham khảo",
 "Ngon",
 "nhạy",
 "nhẹ",
 "ok",
 "ổn",
 "Ổn Áp",
 "ổn định",
 "Pin khỏe",
 "quá tốt",
 "quá tuyệt",
 "Rất hữu dụng",
 "rất ổn",
 "Rất thích",
 "rất thích",
 "rất tốt",
 "sắc nét",
 "Sài ok phết",
 "siêu nhanh",
 "sóng khỏe",
 "sử dụng tốt",
 "suất sắc",
 "tạm ổn",
 "thoải mái",
 "thời gian bảo hành tốt",
 "tiện dụng",
 "tốt",
 "tuyệt vời",
 "Xài êm",
 "xài phê",
 "xịn xò",
 "xử lí nhanh",
 "xử lý tốt",
 "xứng đáng",
 "khỏi bàn",
 " không bị rè",
 "rất hay",
 "rõ ràng",
 " to",
 "tiện lợi",
 "sướng thật",
 "thông minh",
 " ăn đứt",
 "gọn nhẹ",
 "hay",
 "ko giật lag",
 "ko lo hỏng",
 "phong phú"
)

for ($i = 0; $i -lt $totWords.Length; $i++) {
    $r = $i + 1
    $wsTot.Cells.Item($r, 1).Value = $totWords[$i]
}

# --- Populate "Không tốt" sheet (negative sentiment words, insertion order) ---
$khongTotWords = @(
 "ọp ẹp",
 "chán",
 "load chậm",
 "không nghe",
 " không hài lòng ",
 "không được sắc nét",
 "không hiểu",
 "bực mình",
 "độ trễ cao",
 "kém",
 "hư",
 "chậm",
 "sai lầm",
 "yêu",
 "dởm",
 "tẩy chay",
 "nhanh hết pin",
 "thất vọng",
 "rất chán",
 "hơi nặng",
 "hơi nặng",
 " giật lag",
 "đơ",
 "Máy lỗi",
 "rất khó chịu",
 " tụt ghê",
 "rất kém",
 "yếu",
 "kinh khủng",
 "xấu",
 " lỏng lẻo",
 " không ổn định",
 "rất nóng",
 " rất tệ",
 "nghe ko được",
 "mất uy tín",
 "Sai lầm",
 "rè rè",
 "giao diện rố",
 "không nên mua",
 "cực kỳ chậm",
 "sập nguồn",
 "mãi tệ",
 "đơ lác",
 " thấy ghê",
 "quá chậm",
 "quá ì ạch",
 "dễ bị liệt",
 " ức chế",
 "không cứng cáp",
 "đơ lag"
)

for ($i = 0; $i -lt $khongTotWords.Length; $i++) {
    $r = $i + 1
    $wsKhongTot.Cells.Item($r, 1).Value = $khongTotWords[$i]
}

# --- Restore the active selection/view on Sheet1 ---
$ws1.Activate()
$ws1.Range("A300").Select()
